# Grammar pass over the footer/header: bump the "Footer"/"Header" paragraph
# styles (and their linked "Footer Char"/"Header Char" character styles) to
# Arial 11pt, matching what was picked for the rest of the document before
# exporting to PDF.

$d = $word.ActiveDocument

$footer = $d.Styles("Footer")
$footer.Font.Name = "Arial"
$footer.Font.Size = 11

$footerChar = $d.Styles("FooterChar")
$footerChar.Font.Name = "Arial"
$footerChar.Font.Size = 11
$footerChar.Font.SizeBi = 12

$header = $d.Styles("Header")
$header.Font.Name = "Arial"
$header.Font.Size = 11

$headerChar = $d.Styles("HeaderChar")
$headerChar.Font.Name = "Arial"
$headerChar.Font.Size = 11
$headerChar.Font.SizeBi = 12
